$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.624.49"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").Value = "2.119.90"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").Value = "  +0.59%  "
$ws.Range("D5").Value = "'350.17"
$ws.Range("E5").Value = "  +4.69%  "
$ws.Range("D6").Value = "'1.009"
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("D7").Value = "'0.5263"
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("D8").Value = "'0.4514"
$ws.Range("E8").Value = "  -0.91%  "
$ws.Range("D9").Value = "'54.28"
$ws.Range("E9").Value = "  +1.69%  "
$ws.Range("D10").Value = "'0.09086"
$ws.Range("E10").Value = "  +1.86%  "
$ws.Range("D11").Value = "'1.179"
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").Value = "'24.58"
$ws.Range("E12").Value = "  +0.68%  "
$ws.Range("D13").Value = "2.104.78"
$ws.Range("E13").Value = "  +0.48%  "
$ws.Range("D14").Value = "'6.844"
$ws.Range("D15").Value = "'8.085"
$ws.Range("E15").Value = "  +0.80%  "
$ws.Range("D16").Value = "'102.59"
$ws.Range("E16").Value = "  +6.07%  "
$ws.Range("D17").Value = "'0.00001176"
$ws.Range("E17").Value = "  +3.48%  "
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("D19").Value = "'0.06726"
$ws.Range("E19").Value = "  +1.16%  "
$ws.Range("D20").Value = "'19.48"
$ws.Range("E20").Value = "  +1.01%  "
$ws.Range("D21").Value = "'1.010"
$ws.Range("E21").Value = "  +0.68%  "
$ws.Range("D22").Value = "'6.322"
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("D23").Value = "30.674.19"
$ws.Range("E23").Value = "  +0.64%  "
$ws.Range("E24").Value = "  +2.98%  "
$ws.Range("D25").Value = "'2.385"
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("D26").Value = "2.351.35"
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "'22.51"
$ws.Range("E27").Value = "  +0.54%  "
$ws.Range("D28").Value = "'165.25"
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("D29").Value = "'2.559"
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").Value = "'136.61"
$ws.Range("E30").Value = "  +2.77%  "
$ws.Range("D31").Value = "'1.196"
$ws.Range("E31").Value = "  -3.61%  "
$ws.Range("D32").Value = "'0.1077"
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("D33").Value = "'1.663"
$ws.Range("E33").Value = "  -2.44%  "
$ws.Range("D34").Value = "'6.387"
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "'4.017"
$ws.Range("E35").Value = "  +2.12%  "
$ws.Range("D36").Value = "'10.43"
$ws.Range("E36").Value = "  -0.89%  "
$ws.Range("D37").Value = "'5.922"
$ws.Range("E37").Value = "  +5.24%  "
$ws.Range("D38").Value = "'0.02647"
$ws.Range("E38").Value = "  +2.23%  "
$ws.Range("D39").Value = "'0.06873"
$ws.Range("E39").Value = "  +0.54%  "
$ws.Range("D40").Value = "'0.2322"
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("D41").Value = "'12.60"
$ws.Range("E41").Value = "  -1.34%  "
$ws.Range("D42").Value = "'0.6913"
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").Value = "'1.273"
$ws.Range("E43").Value = "  +1.94%  "
$ws.Range("D44").Value = "'14.80"
$ws.Range("E44").Value = "  +4.78%  "
$ws.Range("D45").Value = "'2.335"
$ws.Range("E45").Value = "  -1.04%  "
$ws.Range("D46").Value = "'0.6468"
$ws.Range("E46").Value = "  +1.01%  "
$ws.Range("D47").Value = "'3.751"
$ws.Range("E47").Value = "  +2.43%  "
$ws.Range("D48").Value = "'0.00000000366"
$ws.Range("E48").Value = "  +4.57%  "
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("D50").Value = "'0.07312"
$ws.Range("E50").Value = "  +2.31%  "
$ws.Range("D51").Value = "'82.67"
$ws.Range("E51").Value = "  -0.97%  "
